# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2

# Row 5
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.38
$ws.Range("L5").Value = 6
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.3
$ws.Range("R5").Value = 1.6
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("Z5").Value = 12
$ws.Range("AL5").Value = 51
$ws.Range("AO5").Value = 9

# Row 6
$ws.Range("G6").Value = 1.42
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 1.88
$ws.Range("K6").Value = 2.18
$ws.Range("L6").Value = 7.7
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 6.75
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 2.82
$ws.Range("Q6").Value = 1.98
$ws.Range("U6").Value = 2.15
$ws.Range("V6").Value = 1.55
$ws.Range("W6").Value = 5.2
$ws.Range("X6").Value = 5.7
$ws.Range("Z6").Value = 9
$ws.Range("AA6").Value = 13
$ws.Range("AC6").Value = 8.25
$ws.Range("AD6").Value = 7.6
$ws.Range("AE6").Value = 22
$ws.Range("AH6").Value = 19.5
$ws.Range("AI6").Value = 65
$ws.Range("AJ6").Value = 28
$ws.Range("AK6").Value = 300
$ws.Range("AM6").Value = 110
$ws.Range("AN6").Value = 3.05
$ws.Range("AO6").Value = 6.4
$ws.Range("AP6").Value = 17
$ws.Range("AQ6").Value = 19
$ws.Range("AR6").Value = 50
$ws.Range("AT6").Value = 2.55
$ws.Range("AW6").Value = 9.5
$ws.Range("AX6").Value = 55
$ws.Range("AZ6").Value = 450

# Row 7
$ws.Range("N7").Value = 7.9
